$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Replace the rich-text footnote marker in C3 with plain text, matching the
#    look of the normal data cells (bold Arial Narrow 9, centered, thin border).
$c3 = $ws.Range("C3")
$c3.Value = "31.10. - 02.11."
$c3.HorizontalAlignment = -4108
$c3.Font.Name = "Arial Narrow"
$c3.Font.Bold = $true
$c3.Font.Size = 9

# 2. Remove the footnote legend rows beneath the table (rows 19-22).
$ws.Range("A19:B22").Clear()

# 3. Update the view: zoom level, scroll position and selection.
$ws.Application.ActiveWindow.Zoom = 160
$ws.Range("A19:G23").Select()
$ws.Application.ActiveWindow.ScrollRow = 2
